$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 5.5224
    13 = 6.349999999999997
    16 = 5.103099999999999
    18 = 6.413599999999994
    20 = 8.799099999999999
    26 = 4.830700000000004
    27 = 6.402000000000005
    29 = 5.234399999999999
    35 = 8.703600000000002
    36 = 9.295500000000009
    45 = 5.014800000000005
    55 = 6.616499999999993
    57 = 5.164799999999994
    69 = 5.283499999999995
    76 = 5.0296
    78 = 10.20190000000001
    82 = 5.366400000000001
    83 = 5.210199999999998
    93 = 5.2182
    97 = 5.717099999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
